$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# --- Insert a new row above row 3 (pushes the existing rows 3..51 down to 4..52) ---
$ws.Rows.Item(3).Insert()

# Copy the cell formatting (borders/fill/font/alignment) from the row just below
# (row 4, which used to be row 3 before the insert) onto the freshly inserted row 3,
# so the new row keeps the same look as the rest of the table.
$ws.Range("A4:G4").Copy()
$ws.Range("A3:G3").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Fill in the new row with the "42. Trapping Rain Water" note ---
$ws.Range("A3").Value = "42. Trapping Rain Water"
$ws.Range("B3").Value = "Hard"
$ws.Range("C3").Value = "Array"
$ws.Range("D3").Value = "- Give an array of height representing an elevation map`n- Each bar width is 1`n- return the water can trap after rainning"
$ws.Range("E3").Value = "- We must know at a point which wether this point can be a hole. Mean that this point lower than another ahead bar and an other behind bar`n- We can found the max and min until a point by use 2 array to store these information`n- then at the end, we can use the formular a point minus that min of(max value behind, max value ahead) then it will be the answer`n--> O(n)`n- Because all we need is current index can be a hole or not`n- We can keep two pointer left and right.`n- Firstly, compare these two value`n- Then we can find max of that part, exp maxLeft or maxRight`n- Then just add to answer the amount of current value and max of part`n- Why this algorithms can be used?"
$ws.Range("F3").Value = "- Keep practice,  because I'm kind of not good at algorimths :("

# Row height grew to fit the new (taller) content
$ws.Rows.Item(3).RowHeight = 262.5

# --- Keep the "Hard/Medium/Easy" conditional formatting in sync with the row shift ---
# The block that used to cover A2:G4 must now cover A2:G5 (it grew by the inserted row).
$fcs = $ws.Cells.FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    if ($fc.AppliesTo.Address() -eq "$2:$G$4" -or $fc.AppliesTo.Address() -eq "`$A`$2:`$G`$4") {
        $fc.ModifyAppliesToRange($ws.Range("A2:G5"))
    }
}

# --- View state: select E1 (matches where the note table now starts) ---
$ws.Range("E1").Select()
